# Update "想去人数" (F column) figures across the workbook's sheets to
# reflect the latest generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 8263
    4  = 1937
    5  = 6536
    6  = 163
    7  = 2099
    8  = 584
    9  = 57
    10 = 22
    11 = 59
    15 = 9
    16 = 8646
    17 = 164
    18 = 65
    19 = 196
    26 = 53
    28 = 191
    30 = 16
    31 = 19
    32 = 14
    33 = 2142
    34 = 857
    35 = 500
    37 = 4
    39 = 216
    40 = 157
    42 = 56
    44 = 53
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 27

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$sheet3Updates = @{
    2 = 2332
    4 = 313
}
foreach ($row in $sheet3Updates.Keys) {
    $ws3.Range("F$row").Value = $sheet3Updates[$row]
}

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 2332
    6  = 8263
    8  = 313
    9  = 1937
    10 = 6536
    11 = 2099
    13 = 584
    14 = 22
    17 = 59
    20 = 8646
    21 = 65
    22 = 196
    28 = 53
    30 = 191
    31 = 16
    32 = 19
    33 = 14
    34 = 2142
    35 = 857
    37 = 500
    40 = 216
    41 = 157
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
